$d = $word.ActiveDocument

# 1. Split the 'Norma de recuperação' formula text into separate lines using line breaks
$found1 = $d.Content.Find.Execute("A nota final (MF) do aluno que realizou provas de recuperação dependerá da média do semestre (MS) e da média das provas de recuperação (MR), como segue:MF=5 se 5 ≤MR ≤ (10 - MS); MF = (MS + MR) / 2 se MR > (10 – MS)MF = MS se MR < 5.", $true, $false, $false, $false, $false, $false, 1, $false, "A nota final (MF) do aluno que realizou provas de recuperação dependerá da média do semestre (MS) e da média das provas de recuperação (MR), como segue:^lMF=5 se 5 ≤MR ≤ (10 - MS); ^lMF = (MS + MR) / 2 se MR > (10 – MS)^lMF = MS se MR < 5.", 2)
if (-not $found1) { throw "Find/Replace failed for found1" }

# 2. Split the Bibliografia run into one line break per reference entry
$found2 = $d.Content.Find.Execute("BORGES, A.C. Topografia Aplicada a Engenharia Civil. Vol. I e II. Ed. Edgard Blücher. São Paulo, 1977.BREACH,M.; SCHOFIELD, W. Enginnering Surveying. London: Elsevier, 2007.CALIJURI, M. C.; CUNHA, D. G. F. Engenharia ambiental: Conceitos, tecnologia e gestão. Elsevier, 2012. DUARTE, P.A. Fundamentos de Cartografia. UFSC. Florianópolis, 2006.ESPARTEL,L. Curso de Topografia. Porto Alegre: Globo, 1980FITZ, P. R. Cartografia Básica. São Paulo. Oficina de Textos, 2008. GARCIA, G. J. Sensoriamento Remoto: princípios e interpretação de imagens. São Paulo, Nobel, 1982.JOLY, F. La Cartografia. Barcelona, Ariel, 1982.LIBAUT, A. Geocartografia. Ed. Nacional, Universidade de São Paulo, 1975. LOCH, R. E. N. Cartografia: representação, comunicação e visualização de dados espaciais. Florianópolis: Editora da UFSC, 2006.MARTINELLI, M. Mapas da Geografia e Cartografia Temática. Contexto, 2003.MENEZES, P.M.L; FERNANDES, M.C. Roteiro de cartografia. Oficina de Textos, 2013. MIRANDA,J. Fundamentos de Sistemas de Informações Geográficas. Brasília: Embrapa, 2005.MONICO, J. Posicionamento pelo Navstar-GPS. São Paulo: UNESP, 2000.MOREIRA, M.A. Fundamentos do Sensoreamento Remoto. Editora do INPE, 2001.OLIVEIRA, C. Dicionário cartográfico. Rio de Janeiro, IBGE, 1983.ROBINSON, A. H.; MORRISON, J. L.; MUEHRCKE, P.C.A.; KIMERLING, J; GUPTIL, S. C. Elements of cartography, 6. ed. New York: Wiley, 1995.SILVA, A. Sistemas de informações geo-referenciadas. Campinas: UNICAMP, 2003.VENTURI, L.A.B (Org.) Praticando Geografia: técnicas de campo e laboratório. São Paulo: Oficina dos Textos, 2005. ZUQUETTE, L.V., GANDOLFI, N. Cartografia Geotécnica. Oficina de Textos. São Paulo, 2004.ZUQUETTE, L. V. Geotecnia ambiental. Elsevier, 2015.", $true, $false, $false, $false, $false, $false, 1, $false, "BORGES, A.C. Topografia Aplicada a Engenharia Civil. Vol. I e II. Ed. Edgard Blücher. São Paulo, 1977.^lBREACH,M.; SCHOFIELD, W. Enginnering Surveying. London: Elsevier, 2007.^lCALIJURI, M. C.; CUNHA, D. G. F. Engenharia ambiental: Conceitos, tecnologia e gestão. ^lElsevier, 2012. ^lDUARTE, P.A. Fundamentos de Cartografia. UFSC. Florianópolis, 2006.^lESPARTEL,L. Curso de Topografia. Porto Alegre: Globo, 1980^lFITZ, P. R. Cartografia Básica. São Paulo. Oficina de Textos, 2008. ^lGARCIA, G. J. Sensoriamento Remoto: princípios e interpretação de imagens. São Paulo, Nobel, 1982.^lJOLY, F. La Cartografia. Barcelona, Ariel, 1982.^lLIBAUT, A. Geocartografia. Ed. Nacional, Universidade de São Paulo, 1975. ^lLOCH, R. E. N. Cartografia: representação, comunicação e visualização de dados espaciais. Florianópolis: Editora da UFSC, 2006.^lMARTINELLI, M. Mapas da Geografia e Cartografia Temática. Contexto, 2003.^lMENEZES, P.M.L; FERNANDES, M.C. Roteiro de cartografia. Oficina de Textos, 2013. ^lMIRANDA,J. Fundamentos de Sistemas de Informações Geográficas. Brasília: Embrapa, 2005.^lMONICO, J. Posicionamento pelo Navstar-GPS. São Paulo: UNESP, 2000.^lMOREIRA, M.A. Fundamentos do Sensoreamento Remoto. Editora do INPE, 2001.^lOLIVEIRA, C. Dicionário cartográfico. Rio de Janeiro, IBGE, 1983.^lROBINSON, A. H.; MORRISON, J. L.; MUEHRCKE, P.C.A.; KIMERLING, J; GUPTIL, S. C. Elements of cartography, 6. ed. New York: Wiley, 1995.^lSILVA, A. Sistemas de informações geo-referenciadas. Campinas: UNICAMP, 2003.^lVENTURI, L.A.B (Org.) Praticando Geografia: técnicas de campo e laboratório. São Paulo: Oficina dos Textos, 2005. ^lZUQUETTE, L.V., GANDOLFI, N. Cartografia Geotécnica. Oficina de Textos. São Paulo, 2004.^lZUQUETTE, L. V. Geotecnia ambiental. Elsevier, 2015.", 2)
if (-not $found2) { throw "Find/Replace failed for found2" }

Write-Host "done: $found1 $found2"
